# Applies the cryptos.xlsx data refresh (price / volume updates, row 7/8 coin swap)
# as described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (matches the inlineStr cells in the workbook)
# without leaving a stray NumberFormat behind on the cell style.
foreach ($pair in @(
    @('D2', '64.312.27'),
    @('E2', '  -3.11%  '),
    @('D3', '3.178.58'),
    @('E3', '  -8.01%  '),
    @('E4', '  +0.02%  '),
    @('D5', '565.37'),
    @('E5', '  -3.93%  '),
    @('D6', '169.98'),
    @('E6', '  -3.51%  '),
    @('B7', 'USDC'),
    @('C7', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'),
    @('D7', '1.00'),
    @('E7', '  +0.03%  '),
    @('B8', 'XRP'),
    @('C8', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'),
    @('D8', '0.607'),
    @('E8', '  -1.12%  '),
    @('D9', '3.175.52'),
    @('E9', '  -8.08%  '),
    @('E10', '  -6.33%  '),
    @('D11', '6.65'),
    @('E11', '  -4.36%  '),
    @('E12', '  -4.72%  '),
    @('D13', '3.729.55'),
    @('E13', '  -8.04%  '),
    @('E14', '  +1.59%  '),
    @('D15', '27.43'),
    @('E15', '  -9.22%  '),
    @('D16', '64.326.73'),
    @('E16', '  -2.94%  '),
    @('E17', '  -5.09%  '),
    @('D18', '3.180.59'),
    @('E18', '  -8.02%  '),
    @('D19', '5.74'),
    @('E19', '  -3.73%  '),
    @('D20', '13.02'),
    @('E20', '  -5.74%  '),
    @('D21', '353.16'),
    @('E21', '  -5.47%  '),
    @('D22', '7.20'),
    @('E22', '  -5.49%  '),
    @('E23', '  -0.05%  '),
    @('D24', '69.09'),
    @('E24', '  -5.68%  '),
    @('E25', '  -4.20%  '),
    @('D26', '0.504'),
    @('E26', '  -5.72%  '),
    @('D27', '9.53'),
    @('E27', '  -3.68%  '),
    @('E28', '  -0.86%  '),
    @('E29', '  +0.01%  '),
    @('D30', '5.60'),
    @('E30', '  -4.75%  '),
    @('D32', '1.90'),
    @('E32', '  -5.06%  '),
    @('D33', '22.06'),
    @('E33', '  -6.91%  '),
    @('D34', '6.64'),
    @('E34', '  -5.83%  '),
    @('E35', '  -5.38%  '),
    @('E36', '  -7.03%  '),
    @('D37', '154.47'),
    @('E37', '  -3.70%  '),
    @('D38', '0.819'),
    @('E38', '  -7.48%  '),
    @('D39', '25.92'),
    @('E39', '  -8.53%  '),
    @('E40', '  -1.95%  '),
    @('E41', '  -5.62%  '),
    @('D42', '2.627.65'),
    @('E42', '  -5.04%  '),
    @('D43', '4.18'),
    @('E43', '  -7.17%  '),
    @('E44', '  -6.69%  '),
    @('D45', '39.55'),
    @('E45', '  -1.24%  '),
    @('D46', '0.0656'),
    @('E46', '  -5.28%  '),
    @('D47', '23.76'),
    @('E47', '  -5.81%  '),
    @('D48', '323.87'),
    @('E48', '  -4.42%  '),
    @('E49', '  -7.39%  '),
    @('E50', '  -1.86%  '),
    @('D51', '0.999'),
    @('E51', '  -0.03%  ')
)) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.ClearFormats()
}
